$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 111643.664
$ws.Range("J88").Value = 640.5
$ws.Range("L88").Value = 640.5
$ws.Range("N88").Value = -1452.5
$ws.Range("H91").Value = 111643.664
$ws.Range("J91").Value = 640.5
$ws.Range("L91").Value = 640.5
$ws.Range("N91").Value = -3448.5
$ws.Range("H132").Value = 1615.7646
$ws.Range("I132").Value = 1216.5625
$ws.Range("K132").Value = 3649.6875
$ws.Range("M132").Value = -1119.6875
$ws.Range("H135").Value = 773.7308
$ws.Range("I135").Value = 527.5217
$ws.Range("K135").Value = 4747.6953
$ws.Range("M135").Value = -2212.6953
$ws.Range("H138").Value = 2710.52
$ws.Range("I138").Value = 2058.2727
$ws.Range("J138").Value = 3223
$ws.Range("K138").Value = 6174.8181
$ws.Range("L138").Value = 9669
$ws.Range("M138").Value = -1034.8181
$ws.Range("N138").Value = -19949
$ws.Range("H141").Value = 5004.9165
$ws.Range("I141").Value = 4505.8
$ws.Range("K141").Value = 13517.4
$ws.Range("M141").Value = -8337.400000000001

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3801.38
$ws.Range("I32").Value = 2896.2422
$ws.Range("K32").Value = 2896.2422
$ws.Range("M32").Value = -2609.2422
$ws.Range("H121").Value = 69694
$ws.Range("J121").Value = 69694
$ws.Range("L121").Value = 69694
$ws.Range("N121").Value = -73188
$ws.Range("H124").Value = 44569.125
$ws.Range("J124").Value = 44569.125
$ws.Range("L124").Value = 44569.125
$ws.Range("N124").Value = -54389.125

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 3504.2856
$ws.Range("J22").Value = 6957.5
$ws.Range("L22").Value = 6957.5
$ws.Range("N22").Value = -7303.5
$ws.Range("H134").Value = 2189.9062
$ws.Range("I134").Value = 1836.2593
$ws.Range("J134").Value = 4099.6
$ws.Range("K134").Value = 5508.7779
$ws.Range("L134").Value = 12298.8
$ws.Range("M134").Value = -2973.7779
$ws.Range("N134").Value = -17368.8

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4278.0654
$ws.Range("I31").Value = 2350.4443
$ws.Range("K31").Value = 2350.4443
$ws.Range("M31").Value = -2055.4443
$ws.Range("H34").Value = 4278.0654
$ws.Range("I34").Value = 2350.4443
$ws.Range("K34").Value = 2350.4443
$ws.Range("M34").Value = -2148.4443
$ws.Range("H76").Value = 5560553.5
$ws.Range("I76").Value = 5560553.5
$ws.Range("K76").Value = 5560553.5
$ws.Range("M76").Value = -5560238.5
$ws.Range("H79").Value = 5560553.5
$ws.Range("I79").Value = 5560553.5
$ws.Range("K79").Value = 5560553.5
$ws.Range("M79").Value = -5559461.5
$ws.Range("H94").Value = 3358.6
$ws.Range("I94").Value = 3698
$ws.Range("K94").Value = 3698
$ws.Range("M94").Value = -3247
$ws.Range("H99").Value = 1955622.1
$ws.Range("I99").Value = 2122.5
$ws.Range("J99").Value = 2234693.5
$ws.Range("K99").Value = 2122.5
$ws.Range("L99").Value = 2234693.5
$ws.Range("M99").Value = -624.5
$ws.Range("N99").Value = -2237689.5
$ws.Range("H105").Value = 103045.91
$ws.Range("I105").Value = 140532.12
$ws.Range("J105").Value = 3082.6667
$ws.Range("K105").Value = 140532.12
$ws.Range("L105").Value = 3082.6667
$ws.Range("M105").Value = -138785.12
$ws.Range("N105").Value = -6576.6667
$ws.Range("H109").Value = 43857.145
$ws.Range("J109").Value = 43857.145
$ws.Range("L109").Value = 43857.145
$ws.Range("N109").Value = -45937.145
$ws.Range("H126").Value = 1955622.1
$ws.Range("I126").Value = 2122.5
$ws.Range("J126").Value = 2234693.5
$ws.Range("K126").Value = 6367.5
$ws.Range("L126").Value = 6704080.5
$ws.Range("M126").Value = -3897.5
$ws.Range("N126").Value = -6709020.5
$ws.Range("H132").Value = 1570.2941
$ws.Range("I132").Value = 1199.7241
$ws.Range("J132").Value = 3719.6
$ws.Range("K132").Value = 3599.1723
$ws.Range("L132").Value = 11158.8
$ws.Range("M132").Value = -1069.1723
$ws.Range("N132").Value = -16218.8

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1357
$ws.Range("I5").Value = 1014.5
$ws.Range("K5").Value = 3043.5
$ws.Range("M5").Value = -2931.5
$ws.Range("H14").Value = 199.66667
$ws.Range("I14").Value = 199.66667
$ws.Range("K14").Value = 599.00001
$ws.Range("M14").Value = -426.00001
$ws.Range("H34").Value = 1532.2354
$ws.Range("J34").Value = 3966.6667
$ws.Range("L34").Value = 11900.0001
$ws.Range("N34").Value = -12068.0001
$ws.Range("H37").Value = 47670.855
$ws.Range("J37").Value = 47670.855
$ws.Range("L37").Value = 143012.565
$ws.Range("N37").Value = -143236.565
$ws.Range("H131").Value = 45379.87
$ws.Range("I131").Value = 84396.5
$ws.Range("J131").Value = 2816.2727
$ws.Range("K131").Value = 253189.5
$ws.Range("L131").Value = 8448.8181
$ws.Range("M131").Value = -248149.5
$ws.Range("N131").Value = -18528.8181
$ws.Range("H135").Value = 1357
$ws.Range("I135").Value = 1014.5
$ws.Range("K135").Value = 9130.5
$ws.Range("M135").Value = -6595.5
$ws.Range("H141").Value = 4075.8
$ws.Range("I141").Value = 3417.5557
$ws.Range("J141").Value = 10000
$ws.Range("K141").Value = 10252.6671
$ws.Range("L141").Value = 30000
$ws.Range("M141").Value = -5072.667099999999
$ws.Range("N141").Value = -40360

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 50000
$ws.Range("J32").Value = 50000
$ws.Range("L32").Value = 50000
$ws.Range("N32").Value = -50592
$ws.Range("H42").Value = 48217.5
$ws.Range("J42").Value = 48217.5
$ws.Range("L42").Value = 48217.5
$ws.Range("N42").Value = -49187.5
$ws.Range("H115").Value = 48217.5
$ws.Range("J115").Value = 48217.5
$ws.Range("L115").Value = 48217.5
$ws.Range("N115").Value = -50567.5
$ws.Range("H122").Value = 131391.95
$ws.Range("I122").Value = 141965.44
$ws.Range("K122").Value = 425896.32
$ws.Range("M122").Value = -423446.32
$ws.Range("H132").Value = 2359.889
$ws.Range("I132").Value = 2132.5908
$ws.Range("K132").Value = 6397.7724
$ws.Range("M132").Value = -3867.7724

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 27854.723
$ws.Range("I7").Value = 45979
$ws.Range("J7").Value = 13355.3
$ws.Range("K7").Value = 45979
$ws.Range("L7").Value = 13355.3
$ws.Range("M7").Value = -45867
$ws.Range("N7").Value = -13579.3
$ws.Range("H40").Value = 3475200
$ws.Range("I40").Value = 2866.818
$ws.Range("K40").Value = 2866.818
$ws.Range("M40").Value = -2730.818
$ws.Range("H100").Value = 22179.334
$ws.Range("I100").Value = 31920.1
$ws.Range("J100").Value = 2697.8
$ws.Range("K100").Value = 31920.1
$ws.Range("L100").Value = 2697.8
$ws.Range("M100").Value = -31379.1
$ws.Range("N100").Value = -3779.8
$ws.Range("H122").Value = 20008262
$ws.Range("I122").Value = 8668.286
$ws.Range("K122").Value = 26004.858
$ws.Range("M122").Value = -23554.858
$ws.Range("H126").Value = 27854.723
$ws.Range("I126").Value = 45979
$ws.Range("J126").Value = 13355.3
$ws.Range("K126").Value = 137937
$ws.Range("L126").Value = 40065.89999999999
$ws.Range("M126").Value = -135467
$ws.Range("N126").Value = -45005.89999999999
$ws.Range("H132").Value = 2924.5356
$ws.Range("I132").Value = 2434.75
$ws.Range("K132").Value = 7304.25
$ws.Range("M132").Value = -4774.25
$ws.Range("H135").Value = 68499.664
$ws.Range("J135").Value = 68499.664
$ws.Range("L135").Value = 68499.664
$ws.Range("N135").Value = -78639.664
$ws.Range("H136").Value = 7152.96
$ws.Range("I136").Value = 7556.7334
$ws.Range("J136").Value = 6547.3
$ws.Range("K136").Value = 22670.2002
$ws.Range("L136").Value = 19641.9
$ws.Range("M136").Value = -20120.2002
$ws.Range("N136").Value = -24741.9
$ws.Range("H139").Value = 79606.42999999999
$ws.Range("I139").Value = 66333.336
$ws.Range("J139").Value = 89561.25
$ws.Range("K139").Value = 66333.336
$ws.Range("L139").Value = 89561.25
$ws.Range("M139").Value = -61193.336
$ws.Range("N139").Value = -99841.25

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 50507.855
$ws.Range("J119").Value = 50507.855
$ws.Range("L119").Value = 50507.855
$ws.Range("N119").Value = -60183.855
$ws.Range("H132").Value = 1404893.2
$ws.Range("I132").Value = 2176.3333
$ws.Range("K132").Value = 6528.999899999999
$ws.Range("M132").Value = -3998.999899999999
